$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NavigateURL value for "clients" row from "/clients" to "clients"
$ws.Range("B1").Value = "clients"

# Re-select cell B2 to match the saved selection state
$ws.Range("B2").Select()
